# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on several
# leve rows across the ALC, ARM, BSM, CRP, CUL and LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1446.5714
$ws.Range("I20").Value = 1446.5714
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1446.5714
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1216.5714

$ws.Range("H35").Value = 1446.5714
$ws.Range("I35").Value = 1446.5714
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1446.5714
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1067.5714

$ws.Range("H40").Value = 1813.0541
$ws.Range("I40").Value = 1437.96
$ws.Range("J40").Value = 2594.5
$ws.Range("K40").Value = 1437.96
$ws.Range("L40").Value = 2594.5
$ws.Range("M40").Value = -1262.96
$ws.Range("N40").Value = -2944.5

$ws.Range("H43").Value = 689.7917
$ws.Range("I43").Value = 579.1
$ws.Range("J43").Value = 768.8570999999999
$ws.Range("K43").Value = 579.1
$ws.Range("L43").Value = 768.8570999999999
$ws.Range("M43").Value = -510.1
$ws.Range("N43").Value = -906.8570999999999

$ws.Range("H113").Value = 3292.25
$ws.Range("I113").Value = 2788
$ws.Range("J113").Value = 4132.6665
$ws.Range("K113").Value = 2788
$ws.Range("L113").Value = 4132.6665
$ws.Range("M113").Value = 466
$ws.Range("N113").Value = -10640.6665

$ws.Range("H120").Value = 38150
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 38150
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 38150
$ws.Range("N120").Value = -47826

$ws.Range("H129").Value = 1610.4884
$ws.Range("I129").Value = 580.8570999999999
$ws.Range("J129").Value = 2107.5518
$ws.Range("K129").Value = 1742.5713
$ws.Range("L129").Value = 6322.655400000001
$ws.Range("M129").Value = 3257.4287
$ws.Range("N129").Value = -16322.6554

$ws.Range("H132").Value = 5261.754
$ws.Range("I132").Value = 4293.7334
$ws.Range("J132").Value = 7984.3125
$ws.Range("K132").Value = 12881.2002
$ws.Range("L132").Value = 23952.9375
$ws.Range("M132").Value = -10351.2002
$ws.Range("N132").Value = -29012.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 27000
$ws.Range("I57").Value = 27000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 27000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -26516

$ws.Range("H97").Value = 1762.5
$ws.Range("I97").Value = 1820
$ws.Range("J97").Value = 1666.6666
$ws.Range("K97").Value = 1820
$ws.Range("L97").Value = 1666.6666
$ws.Range("M97").Value = -1324
$ws.Range("N97").Value = -2658.6666

$ws.Range("H122").Value = 1224.5
$ws.Range("I122").Value = 882.44446
$ws.Range("J122").Value = 1566.5555
$ws.Range("K122").Value = 2647.33338
$ws.Range("L122").Value = 4699.666499999999
$ws.Range("M122").Value = -197.33338
$ws.Range("N122").Value = -9599.666499999999

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 30000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -27530

$ws.Range("H132").Value = 1816660.8
$ws.Range("I132").Value = 5719.387
$ws.Range("J132").Value = 3294007.5
$ws.Range("K132").Value = 17158.161
$ws.Range("L132").Value = 9882022.5
$ws.Range("M132").Value = -14628.161
$ws.Range("N132").Value = -9887082.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 302.5
$ws.Range("I12").Value = 505
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 505
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = -337
$ws.Range("N12").Value = -436

$ws.Range("H82").Value = 15033.134
$ws.Range("I82").Value = 2521.4
$ws.Range("J82").Value = 40056.6
$ws.Range("K82").Value = 2521.4
$ws.Range("L82").Value = 40056.6
$ws.Range("M82").Value = -2138.4
$ws.Range("N82").Value = -40822.6

$ws.Range("H85").Value = 15033.134
$ws.Range("I85").Value = 2521.4
$ws.Range("J85").Value = 40056.6
$ws.Range("K85").Value = 2521.4
$ws.Range("L85").Value = 40056.6
$ws.Range("M85").Value = -1195.4
$ws.Range("N85").Value = -42708.6

$ws.Range("H94").Value = 1559.4
$ws.Range("I94").Value = 1535.0435
$ws.Range("J94").Value = 1639.4286
$ws.Range("K94").Value = 1535.0435
$ws.Range("L94").Value = 1639.4286
$ws.Range("M94").Value = -1084.0435
$ws.Range("N94").Value = -2541.4286

$ws.Range("H134").Value = 2475.6445
$ws.Range("I134").Value = 1571.84
$ws.Range("J134").Value = 3605.4
$ws.Range("K134").Value = 4715.52
$ws.Range("L134").Value = 10816.2
$ws.Range("M134").Value = -2180.52
$ws.Range("N134").Value = -15886.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 15340.3
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15340.3
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 15340.3
$ws.Range("N43").Value = -15708.3

$ws.Range("H58").Value = 7876.6313
$ws.Range("I58").Value = 4613
$ws.Range("J58").Value = 13471.429
$ws.Range("K58").Value = 4613
$ws.Range("L58").Value = 13471.429
$ws.Range("M58").Value = -4410
$ws.Range("N58").Value = -13877.429

$ws.Range("H62").Value = 19366.666
$ws.Range("I62").Value = 3225
$ws.Range("J62").Value = 51650
$ws.Range("K62").Value = 3225
$ws.Range("L62").Value = 51650
$ws.Range("M62").Value = -2601
$ws.Range("N62").Value = -52898

$ws.Range("H65").Value = 19366.666
$ws.Range("I65").Value = 3225
$ws.Range("J65").Value = 51650
$ws.Range("K65").Value = 16125
$ws.Range("L65").Value = 258250
$ws.Range("M65").Value = -13005
$ws.Range("N65").Value = -264490

$ws.Range("H101").Value = 15340.3
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 15340.3
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 15340.3
$ws.Range("N101").Value = -21830.3

$ws.Range("H136").Value = 7876.6313
$ws.Range("I136").Value = 4613
$ws.Range("J136").Value = 13471.429
$ws.Range("K136").Value = 13839
$ws.Range("L136").Value = 40414.287
$ws.Range("M136").Value = -11289
$ws.Range("N136").Value = -45514.287

$ws.Range("H138").Value = 40382.5
$ws.Range("I138").Value = 20000
$ws.Range("J138").Value = 42235.453
$ws.Range("K138").Value = 20000
$ws.Range("L138").Value = 42235.453
$ws.Range("M138").Value = -14860
$ws.Range("N138").Value = -52515.453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6000
$ws.Range("N34").Value = -6168

$ws.Range("H131").Value = 836.88464
$ws.Range("I131").Value = 468.57144
$ws.Range("J131").Value = 1266.5834
$ws.Range("K131").Value = 1405.71432
$ws.Range("L131").Value = 3799.7502
$ws.Range("M131").Value = 3634.28568
$ws.Range("N131").Value = -13879.7502

$ws.Range("H132").Value = 1087.2812
$ws.Range("I132").Value = 1049.3572
$ws.Range("J132").Value = 1116.7778
$ws.Range("K132").Value = 9444.2148
$ws.Range("L132").Value = 10051.0002
$ws.Range("M132").Value = -6914.2148
$ws.Range("N132").Value = -15111.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 24129.5
$ws.Range("I99").Value = 33259
$ws.Range("J99").Value = 15000
$ws.Range("K99").Value = 33259
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -30264
$ws.Range("N99").Value = -20990

$ws.Range("H122").Value = 9604.294
$ws.Range("I122").Value = 16857.143
$ws.Range("J122").Value = 4527.3
$ws.Range("K122").Value = 50571.429
$ws.Range("L122").Value = 13581.9
$ws.Range("M122").Value = -48121.429
$ws.Range("N122").Value = -18481.9
